$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.1509980484562094, 1.195282167228728, 7.57942700223539, 2.753075916540514, 2.775751383040755, 52)
    3  = @(0.03079452806175957, 1.189334643345628, 5.574265677164387, 2.360988283995578, 2.384278448858261, 51)
    4  = @(-0.101817367501059, 1.016702669087111, 4.454714751352228, 2.11061951837659, 2.129565437855708, 50)
    5  = @(0.05611136314148432, 1.148195161250686, 5.204774752689769, 2.281397543763421, 2.304342309971512, 49)
    6  = @(-0.06628332611628297, 1.009343296419617, 4.545192847316543, 2.131945789019163, 2.153465142626193, 48)
    7  = @(-0.01273090096952459, 1.122250521844307, 5.511138675828411, 2.347581452437468, 2.380847133747632, 36)
    8  = @(-0.01901641356679982, 1.14391309009801, 5.569593947515463, 2.359998717693605, 2.394375341132336, 35)
    9  = @(0.05598467855657813, 1.513431136887689, 9.181581124177706, 3.030112394644414, 3.11742770206278, 18)
    10 = @(-0.7620237227489945, 1.215958780976282, 6.693346821907633, 2.587150328432353, 2.593054826775349, 11)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
